# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Swap two pairs of country names (table re-sorted by "Casos totales") ---
# Rows 38/39: Kazajistan <-> Paises Bajos
$ws.Range("A38").Value = "Paises Bajos"
$ws.Range("A39").Value = "Kazajistan"

# Rows 205/206: Santa Lucia <-> Timor Oriental
$ws.Range("A205").Value = "Timor Oriental"
$ws.Range("A206").Value = "Santa Lucia"

# --- Update the "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 26 de Septiembre de 2020 a las 14:34"

# --- Update numeric data cells ---

# Row 4 (Estados Unidos)
$ws.Range("B4").Value = 7245723
$ws.Range("C4").Value = 1539
$ws.Range("D4").Value = 4481029
$ws.Range("E4").Value = 2556211
$ws.Range("G4").Value = 43
$ws.Range("H4").Value = 208483

# Row 38 (now Paises Bajos)
$ws.Range("B38").Value = 108631
$ws.Range("C38").Value = 2713
$ws.Range("D38").Value = 0
$ws.Range("E38").Value = 0
$ws.Range("G38").Value = 38
$ws.Range("H38").Value = 6366

# Row 39 (now Kazajistan)
$ws.Range("B39").Value = 107659
$ws.Range("C39").Value = 69
$ws.Range("D39").Value = 102530
$ws.Range("E39").Value = 3430
$ws.Range("H39").Value = 1699

# Row 68
$ws.Range("B68").Value = 39895
$ws.Range("C68").Value = 108
$ws.Range("D68").Value = 37523
$ws.Range("E68").Value = 1787
$ws.Range("G68").Value = 2
$ws.Range("H68").Value = 585

# Row 70
$ws.Range("B70").Value = 38253
$ws.Range("C70").Value = 290
$ws.Range("D70").Value = 27704
$ws.Range("E70").Value = 10264
$ws.Range("G70").Value = 7
$ws.Range("H70").Value = 285

# Row 78
$ws.Range("B78").Value = 27016
$ws.Range("C78").Value = 16
$ws.Range("E78").Value = 1575

# Row 79
$ws.Range("B79").Value = 26797
$ws.Range("C79").Value = 233
$ws.Range("D79").Value = 19746
$ws.Range("E79").Value = 6231
$ws.Range("G79").Value = 12
$ws.Range("H79").Value = 820

# Row 80
$ws.Range("B80").Value = 26213
$ws.Range("C80").Value = 619
$ws.Range("D80").Value = 19350
$ws.Range("E80").Value = 6215
$ws.Range("G80").Value = 1
$ws.Range("H80").Value = 648

# Row 88
$ws.Range("B88").Value = 16257
$ws.Range("C88").Value = 36
$ws.Range("D88").Value = 14922
$ws.Range("E88").Value = 1106
$ws.Range("G88").Value = 1
$ws.Range("H88").Value = 229

# Row 126
$ws.Range("E126").Value = 3113
$ws.Range("G126").Value = 1
$ws.Range("H126").Value = 28
